$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos -> now holds the professor string (content bug reproduced from source diff)
$professor = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("B10").Value = $professor
$ws.Range("C10").Value = $professor

# Row 13: Programa resumido: / Semestral
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14: Short syllabus: / short syllabus English text (unchanged text, only label shifts)
$ws.Range("A14").Value = "Short syllabus:"
$shortSyllabusEn = "Review. Mineral composition of the soil solid phase. Composition of the soil organic solid phase. Soil solution. Surface phenomena. Acid soils and salt affected soils. Flooded soils. Soil fertility evaluation methods (Practice). Chemical analyzes of the soil for fertility purposes (Practice)."
$ws.Range("B14").Value = $shortSyllabusEn
$ws.Range("C14").Value = $shortSyllabusEn

# Row 15: Programa: / now holds the activation date string (content bug reproduced from source diff)
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"
$ws.Rows(15).RowHeight = 120

# Row 16: Syllabus: / full English syllabus text
$ws.Range("A16").Value = "Syllabus:"
$syllabusEn = "REVIEW. Soil concepts, profile, composition, characteristics and properties. MINERAL COMPOSITION OF THE SOIL SOLID PHASE. Soil minerals. Major classes of minerals. Origin of electric charges. COMPOSITION OF THE SOIL ORGANIC SOLID PHASE. Composition and structure of soil organic matter, functions and reactions, organic matter and management systems. SOIL SOLUTION. Composition of the soil solution, organic molecules dissolved in the soil solution, concentration and ion activity, obtaining the soil solution. SURFACE PHENOMENA. Origin of electric charges, adsorption descriptive models, zero electric charge point. Field practice class: Soil fertility assessment methods: Land sampling: sampling planning and land sampling. Laboratory Practice Class: Chemical analyzes of the soil for fertility purposes: extractors and analytical methods."
$ws.Range("B16").Value = $syllabusEn
$ws.Range("C16").Value = $syllabusEn

# Row 17: Avaliação: (label only now, B/C cleared)
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows(17).RowHeight = 15

# Row 18: Método: / now holds the professor string again (content bug reproduced from source diff)
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = $professor
$ws.Range("C18").Value = $professor
$ws.Rows(18).RowHeight = 60

# Row 19: Critério: / method-of-evaluation text
$ws.Range("A19").Value = "Critério:"
$metodoTexto = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas notas serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto

# Row 20: Norma de recuperação: / criteria text
$ws.Range("A20").Value = "Norma de recuperação:"
$criterioTexto = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("B20").Value = $criterioTexto
$ws.Range("C20").Value = $criterioTexto

# Row 21: Bibliografia: / recovery-norm text
$ws.Range("A21").Value = "Bibliografia:"
$normaTexto = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("B21").Value = $normaTexto
$ws.Range("C21").Value = $normaTexto
$ws.Rows(21).RowHeight = 120

# Row 22: Requisitos: (label only now, B/C cleared)
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows(22).RowHeight = 15

# Row 23: label cleared, B/C now hold the weak-requisite text
$ws.Range("A23").ClearContents()
$reqTexto = "LOB1206 -  Solos I  (Requisito fraco)`n"
$ws.Range("B23").Value = $reqTexto
$ws.Range("C23").Value = $reqTexto
$ws.Rows(23).RowHeight = 30

# Row 24 no longer exists in the target layout
$ws.Rows(24).Delete()
